# Daily data codeshare - September/October daily actuals update + October N31 ratio + September conditional formatting

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# September sheet: update B2:E31 with revised daily figures
# ---------------------------------------------------------------------------
$wsSep = $wb.Worksheets.Item("September")

$sepArr = New-Object 'object[,]' 30,4
$sepArr[0,0]=251; $sepArr[0,1]=39; $sepArr[0,2]=2; $sepArr[0,3]=210
$sepArr[1,0]=29; $sepArr[1,1]=27; $sepArr[1,2]=0; $sepArr[1,3]=2
$sepArr[2,0]=226; $sepArr[2,1]=22; $sepArr[2,2]=4; $sepArr[2,3]=200
$sepArr[3,0]=243; $sepArr[3,1]=30; $sepArr[3,2]=6; $sepArr[3,3]=207
$sepArr[4,0]=212; $sepArr[4,1]=50; $sepArr[4,2]=5; $sepArr[4,3]=157
$sepArr[5,0]=232; $sepArr[5,1]=30; $sepArr[5,2]=2; $sepArr[5,3]=200
$sepArr[6,0]=191; $sepArr[6,1]=27; $sepArr[6,2]=5; $sepArr[6,3]=159
$sepArr[7,0]=241; $sepArr[7,1]=33; $sepArr[7,2]=6; $sepArr[7,3]=202
$sepArr[8,0]=263; $sepArr[8,1]=38; $sepArr[8,2]=3; $sepArr[8,3]=222
$sepArr[9,0]=229; $sepArr[9,1]=42; $sepArr[9,2]=1; $sepArr[9,3]=186
$sepArr[10,0]=251; $sepArr[10,1]=43; $sepArr[10,2]=4; $sepArr[10,3]=204
$sepArr[11,0]=263; $sepArr[11,1]=46; $sepArr[11,2]=6; $sepArr[11,3]=211
$sepArr[12,0]=250; $sepArr[12,1]=33; $sepArr[12,2]=7; $sepArr[12,3]=210
$sepArr[13,0]=234; $sepArr[13,1]=34; $sepArr[13,2]=4; $sepArr[13,3]=196
$sepArr[14,0]=257; $sepArr[14,1]=28; $sepArr[14,2]=5; $sepArr[14,3]=224
$sepArr[15,0]=261; $sepArr[15,1]=40; $sepArr[15,2]=4; $sepArr[15,3]=217
$sepArr[16,0]=243; $sepArr[16,1]=39; $sepArr[16,2]=3; $sepArr[16,3]=201
$sepArr[17,0]=254; $sepArr[17,1]=40; $sepArr[17,2]=5; $sepArr[17,3]=209
$sepArr[18,0]=242; $sepArr[18,1]=33; $sepArr[18,2]=5; $sepArr[18,3]=204
$sepArr[19,0]=248; $sepArr[19,1]=21; $sepArr[19,2]=11; $sepArr[19,3]=216
$sepArr[20,0]=222; $sepArr[20,1]=37; $sepArr[20,2]=8; $sepArr[20,3]=177
$sepArr[21,0]=237; $sepArr[21,1]=41; $sepArr[21,2]=3; $sepArr[21,3]=193
$sepArr[22,0]=233; $sepArr[22,1]=32; $sepArr[22,2]=2; $sepArr[22,3]=199
$sepArr[23,0]=244; $sepArr[23,1]=39; $sepArr[23,2]=1; $sepArr[23,3]=204
$sepArr[24,0]=247; $sepArr[24,1]=36; $sepArr[24,2]=7; $sepArr[24,3]=204
$sepArr[25,0]=224; $sepArr[25,1]=26; $sepArr[25,2]=8; $sepArr[25,3]=190
$sepArr[26,0]=256; $sepArr[26,1]=34; $sepArr[26,2]=7; $sepArr[26,3]=215
$sepArr[27,0]=222; $sepArr[27,1]=28; $sepArr[27,2]=6; $sepArr[27,3]=188
$sepArr[28,0]=267; $sepArr[28,1]=46; $sepArr[28,2]=7; $sepArr[28,3]=214
$sepArr[29,0]=256; $sepArr[29,1]=34; $sepArr[29,2]=3; $sepArr[29,3]=219
$wsSep.Range("B2:E31").Value = $sepArr

# Rows 5-31 get the existing "red font" style (cellXfs index 2 / red Calibri font)
# applied, matching the same styling already used on the other monthly sheets.
$wsSep.Range("B5:E31").Font.Color = 255

# ---------------------------------------------------------------------------
# October sheet: update B2:E30 with revised daily figures, extend styling,
# add empty styled cell F12, and add N31 ratio formula
# ---------------------------------------------------------------------------
$wsOct = $wb.Worksheets.Item("October")

$octArr = New-Object 'object[,]' 29,4
$octArr[0,0]=245; $octArr[0,1]=40; $octArr[0,2]=2; $octArr[0,3]=203
$octArr[1,0]=233; $octArr[1,1]=28; $octArr[1,2]=1; $octArr[1,3]=204
$octArr[2,0]=228; $octArr[2,1]=3; $octArr[2,2]=3; $octArr[2,3]=192
$octArr[3,0]=255; $octArr[3,1]=34; $octArr[3,2]=5; $octArr[3,3]=216
$octArr[4,0]=210; $octArr[4,1]=32; $octArr[4,2]=5; $octArr[4,3]=173
$octArr[5,0]=256; $octArr[5,1]=33; $octArr[5,2]=7; $octArr[5,3]=216
$octArr[6,0]=245; $octArr[6,1]=32; $octArr[6,2]=6; $octArr[6,3]=207
$octArr[7,0]=219; $octArr[7,1]=38; $octArr[7,2]=5; $octArr[7,3]=176
$octArr[8,0]=229; $octArr[8,1]=28; $octArr[8,2]=5; $octArr[8,3]=196
$octArr[9,0]=239; $octArr[9,1]=36; $octArr[9,2]=3; $octArr[9,3]=200
$octArr[10,0]=240; $octArr[10,1]=33; $octArr[10,2]=4; $octArr[10,3]=203
$octArr[11,0]=209; $octArr[11,1]=31; $octArr[11,2]=2; $octArr[11,3]=176
$octArr[12,0]=251; $octArr[12,1]=30; $octArr[12,2]=9; $octArr[12,3]=212
$octArr[13,0]=246; $octArr[13,1]=34; $octArr[13,2]=9; $octArr[13,3]=203
$octArr[14,0]=239; $octArr[14,1]=37; $octArr[14,2]=3; $octArr[14,3]=199
$octArr[15,0]=235; $octArr[15,1]=32; $octArr[15,2]=4; $octArr[15,3]=199
$octArr[16,0]=243; $octArr[16,1]=30; $octArr[16,2]=9; $octArr[16,3]=204
$octArr[17,0]=254; $octArr[17,1]=40; $octArr[17,2]=2; $octArr[17,3]=212
$octArr[18,0]=226; $octArr[18,1]=38; $octArr[18,2]=5; $octArr[18,3]=183
$octArr[19,0]=228; $octArr[19,1]=45; $octArr[19,2]=6; $octArr[19,3]=177
$octArr[20,0]=250; $octArr[20,1]=31; $octArr[20,2]=4; $octArr[20,3]=215
$octArr[21,0]=222; $octArr[21,1]=18; $octArr[21,2]=7; $octArr[21,3]=197
$octArr[22,0]=220; $octArr[22,1]=35; $octArr[22,2]=8; $octArr[22,3]=177
$octArr[23,0]=233; $octArr[23,1]=27; $octArr[23,2]=7; $octArr[23,3]=199
$octArr[24,0]=252; $octArr[24,1]=29; $octArr[24,2]=6; $octArr[24,3]=217
$octArr[25,0]=203; $octArr[25,1]=34; $octArr[25,2]=3; $octArr[25,3]=166
$octArr[26,0]=230; $octArr[26,1]=13; $octArr[26,2]=28; $octArr[26,3]=189
$octArr[27,0]=246; $octArr[27,1]=17; $octArr[27,2]=21; $octArr[27,3]=208
$octArr[28,0]=202; $octArr[28,1]=15; $octArr[28,2]=24; $octArr[28,3]=163
$wsOct.Range("B2:E30").Value = $octArr

# Apply the red font style across the whole data block (B2:E32), matching rows
# 31:32 which already carried it
$wsOct.Range("B2:E32").Font.Color = 255
$wsOct.Range("F12").Font.Color = 255

$wsOct.Range("N31").Formula = "=466/175"

# ---------------------------------------------------------------------------
# September: add "Equal To $G$14" conditional formatting on H12:H16
# (Highlight Cells Rule -> Equal To..., default Light Red Fill / Dark Red Text)
# ---------------------------------------------------------------------------
$cfRange = $wsSep.Range("H12:H16")
$cf = $cfRange.FormatConditions.Add(1, 3, "=`$G`$14")
$cf.Font.Color = 393372
$cf.Interior.Color = 13551615

# ---------------------------------------------------------------------------
# View state: move the selected/active tab from November to September,
# update the in-sheet selections, and October's scroll position
# ---------------------------------------------------------------------------
$wsOct.Range("H13:J17").Select()

$wsSep.Range("C4").Select()
